# Insert a new weekly price record as row 5 of the "Zapallo" sheet.
# All subsequent rows (previously 5..97) shift down to 6..98, and the
# dimension grows from A1:R97 to A1:R98, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, empty row above current row 5 (row 5 and below shift down).
$ws.Rows.Item(5).Insert()

# Populate the new row 5 with the new weekly observation.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44496
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112045
$ws.Range("G5").Value = "Zapallo"
$ws.Range("H5").Value = "Camote"
$ws.Range("I5").Value = "1a (guarda)"
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = 850
$ws.Range("N5").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("P5").Value = 850
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
